# تعديل تلقائي في شيت Card15 by admin at 2025-11-02 08:02:26
#
# Row 2 of the "Card15" lookup sheet had a leftover/incorrect card number
# ("2") in column A and a truly-empty "Serviced by" cell (O2). Correct the
# card number to match this sheet's card (15) and fill O2 with the same
# "nan" placeholder used by every other data cell in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card15")

# Force column A to stay text (matches the rest of the "card" column,
# which stores its numbers as text) instead of letting Excel infer a
# numeric type for the digit-only string, then drop the number-format
# override so no stray style survives on the cell.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "15"
$ws.Range("A2").ClearFormats()

# O2 was a blank placeholder cell; fill it with "nan" like its neighbors.
$ws.Range("O2").Value = "nan"
